$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 16) with the 2025-05-18 bitcoin buy entry,
# matching the existing rows' layout: Date, Coins, Price, Cost.
$row = 16

# Force column A to be stored as text (matching the other recently
# appended date rows such as A10, A12:A15, which are inline/shared
# strings rather than date serials), then drop the temporary
# number-format override so the cell keeps no explicit style,
# consistent with the existing text-date rows.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "05/18/2025"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = 0.00048095
$ws.Cells.Item($row, 3).Value = 103960.9106975777
$ws.Cells.Item($row, 4).Value = 50
